# Refresh scraped cryptocurrency Price / Volume(1h) figures as published
# by the coinranking.com data pull (GitHub Actions symbol-list update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "275.39" },
    @{ Cell = "E2"; Value = "-1.04%" },
    @{ Cell = "D3"; Value = "26.55" },
    @{ Cell = "E3"; Value = "-2.80%" },
    @{ Cell = "D4"; Value = "4.893" },
    @{ Cell = "E4"; Value = "2.17%" },
    @{ Cell = "D5"; Value = "0.06344" },
    @{ Cell = "E5"; Value = "0.94%" },
    @{ Cell = "D6"; Value = "6.888" },
    @{ Cell = "E6"; Value = "-0.24%" },
    @{ Cell = "D7"; Value = "3.311" },
    @{ Cell = "E7"; Value = "1.39%" },
    @{ Cell = "D8"; Value = "1.285" },
    @{ Cell = "E8"; Value = "36.15%" },
    @{ Cell = "D9"; Value = "0.8680" },
    @{ Cell = "E9"; Value = "-0.99%" },
    @{ Cell = "D10"; Value = "0.1528" },
    @{ Cell = "E10"; Value = "4.82%" },
    @{ Cell = "D11"; Value = "0.05022" },
    @{ Cell = "E11"; Value = "-2.80%" },
    @{ Cell = "D12"; Value = "0.07409" },
    @{ Cell = "E12"; Value = "1.76%" },
    @{ Cell = "D13"; Value = "0.02937" },
    @{ Cell = "E13"; Value = "-5.36%" },
    @{ Cell = "D14"; Value = "0.09041" },
    @{ Cell = "E14"; Value = "-0.33%" },
    @{ Cell = "D15"; Value = "0.001577" },
    @{ Cell = "E15"; Value = "0.90%" },
    @{ Cell = "D16"; Value = "0.0006335" },
    @{ Cell = "E16"; Value = "0.99%" },
    @{ Cell = "D17"; Value = "0.005875" },
    @{ Cell = "E17"; Value = "-0.04%" },
    @{ Cell = "D18"; Value = "3.446" },
    @{ Cell = "E18"; Value = "-0.01%" },
    @{ Cell = "D19"; Value = "2.271" },
    @{ Cell = "E19"; Value = "-0.58%" },
    @{ Cell = "E21"; Value = "0.91%" },
    @{ Cell = "D22"; Value = "3.898" },
    @{ Cell = "E22"; Value = "1.11%" },
    @{ Cell = "D23"; Value = "0.04368" },
    @{ Cell = "E23"; Value = "1.08%" },
    @{ Cell = "D24"; Value = "0.001179" },
    @{ Cell = "D25"; Value = "0.004257" },
    @{ Cell = "E25"; Value = "-0.50%" },
    @{ Cell = "D26"; Value = "0.0001198" },
    @{ Cell = "E26"; Value = "-0.19%" },
    @{ Cell = "D27"; Value = "0.0001676" },
    @{ Cell = "E27"; Value = "-0.90%" },
    @{ Cell = "E40"; Value = "0.75%" },
    @{ Cell = "D41"; Value = "0.006945" },
    @{ Cell = "E41"; Value = "8.06%" },
    @{ Cell = "D42"; Value = "0.1170" },
    @{ Cell = "E42"; Value = "1.10%" },
    @{ Cell = "E43"; Value = "-1.44%" },
    @{ Cell = "D44"; Value = "0.01075" },
    @{ Cell = "E44"; Value = "-9.34%" },
    @{ Cell = "D45"; Value = "0.00005265" },
    @{ Cell = "E45"; Value = "2.62%" },
    @{ Cell = "E46"; Value = "-37.48%" },
    @{ Cell = "D47"; Value = "0.01997" },
    @{ Cell = "E47"; Value = "-11.33%" }
)

foreach ($u in $updates) {
    # Force Text storage before assigning so the numeric-looking
    # price/percentage strings are kept verbatim, matching the
    # worksheet's existing inline-string cell storage.
    $ws.Range($u.Cell).NumberFormat = "@"
    $ws.Range($u.Cell).Value = $u.Value
}
